# ReplyChat.xlsx edit: add "appointment change" rows (23-25), resize row 17,
# and update the sheet selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: shrink the custom row height from 51 to 42 -------------------
$ws.Rows(17).RowHeight = 42

# --- New data rows 23-25 ("appointment change" variants) ------------------
# Row 23
$ws.Range("A23").Value = "Okay USER :) I just need some information in order to make an appointment change. "
$ws.Range("B23").Value = 6
$ws.Range("A23").WrapText = $true
$ws.Rows(23).RowHeight = 18

# Row 24
$ws.Range("A24").Value = "For the appointment change, I just need to confirmation details from you USER :)"
$ws.Range("B24").Value = 6
$ws.Range("A24").WrapText = $true
$ws.Rows(24).RowHeight = 17

# Row 25
$ws.Range("A25").Value = "For the appointment change process, I just need some details from you :)"
$ws.Range("B25").Value = 6
$ws.Range("A25").WrapText = $true
$ws.Rows(25).RowHeight = 17

# --- Match column B number formatting / center alignment (already inherited
# from the column style) and match column A wrap style used by rows 16-22.

# --- Update the visible selection/scroll position --------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G16").Select()
